$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.323.49'
$ws.Range("E2").Value = '  -0.85%  '
$ws.Range("D3").Value = '1.710.96'
$ws.Range("E3").Value = '  -0.69%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.39'
$ws.Range("E5").Value = '  -0.62%  '
$ws.Range("E6").Value = '  -1.41%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06712'
$ws.Range("E8").Value = '  +1.67%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2664'
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.86'
$ws.Range("E10").Value = '  -3.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07677'
$ws.Range("E11").Value = '  -0.61%  '
$ws.Range("D13").Value = '1.946.18'
$ws.Range("E13").Value = '  -0.68%  '
$ws.Range("D14").Value = '1.709.37'
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5820'
$ws.Range("E15").Value = '  -0.35%  '
$ws.Range("D16").Value = '0.0₅8230'
$ws.Range("E16").Value = '  -1.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.20'
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("D18").Value = '27.333.90'
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '226.19'
$ws.Range("E19").Value = '  +2.59%  '
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.625'
$ws.Range("E21").Value = '  -2.13%  '
$ws.Range("E22").Value = '  -2.49%  '
$ws.Range("E23").Value = '  -1.24%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.60'
$ws.Range("E25").Value = '  -2.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.698'
$ws.Range("E26").Value = '  -2.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1207'
$ws.Range("E27").Value = '  -2.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.235'
$ws.Range("E28").Value = '  -2.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '16.35'
$ws.Range("E29").Value = '  -1.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05360'
$ws.Range("E30").Value = '  -3.85%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.291'
$ws.Range("E31").Value = '  -0.76%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.478'
$ws.Range("E32").Value = '  -2.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.430'
$ws.Range("E33").Value = '  -0.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.636'
$ws.Range("E34").Value = '  -1.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.871'
$ws.Range("E35").Value = '  +1.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9498'
$ws.Range("E36").Value = '  -1.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.395'
$ws.Range("E37").Value = '  -1.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5834'
$ws.Range("E38").Value = '  -2.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01636'
$ws.Range("E39").Value = '  -0.82%  '
$ws.Range("D40").Value = '1.081.63'
$ws.Range("E40").Value = '  +2.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.786'
$ws.Range("E41").Value = '  -2.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8399'
$ws.Range("E43").Value = '  -1.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.78'
$ws.Range("E44").Value = '  -0.45%  '
$ws.Range("D45").Value = '1.853.28'
$ws.Range("E45").Value = '  -0.71%  '
$ws.Range("E46").Value = '  +0.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '57.69'
$ws.Range("E47").Value = '  -2.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4533'
$ws.Range("E48").Value = '  +2.25%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.074'
$ws.Range("E49").Value = '  -1.70%  '
$ws.Range("B50").Value = 'Frax'
$ws.Range("C50").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.002'
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("E51").Value = '  -0.24%  '
